# Gatopolis Manager - split "NOME COMPLETO DO ALUNO" into first/last name
# columns (modelo relacional: Name -> first_name / last_name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column A's width so the freshly inserted column matches it.
$colWidthA = $ws.Columns("A").ColumnWidth

# Insert a new column before B; this shifts the former B:G headers one
# column to the right (B:G -> C:H) and keeps their contents/styles intact.
$ws.Columns("B").Insert()
$ws.Columns("B").ColumnWidth = $colWidthA

# Column A used to hold the full name; now it holds the first name only,
# and the newly inserted column B holds the surname.
$ws.Range("A1").Value2 = "NOME  DO ALUNO"
$ws.Range("B1").Value2 = "SOBRENOME DO ALUNO"

# Match the author's final selection/cursor position.
$ws.Range("C6").Select() | Out-Null
